# Applies the TC40 edit: insert two new "scroll" steps after row 27 (TINY_SCROLL_DOWN pattern),
# shifting subsequent rows down by 2, and restyle the header cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 28 so the existing rows 28.. shift down to 30..
$ws.Rows.Item(28).Resize(2).EntireRow.Insert()

# Row 28: SCROLL_WEBELEMENT / Quickorder_SCROLL
$ws.Cells.Item(28, 1).Value = ""
$ws.Cells.Item(28, 2).Value = "SCROLL_WEBELEMENT"
$ws.Cells.Item(28, 3).Value = "Quickorder_SCROLL"
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(28, 5).Value = ""

# Row 29: TINY_SCROLL_DOWN
$ws.Cells.Item(29, 1).Value = ""
$ws.Cells.Item(29, 2).Value = "TINY_SCROLL_DOWN"
$ws.Cells.Item(29, 3).Value = ""
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = ""

# Apply styles matching the surrounding pattern (same as rows 30/31 after insertion)
$ws.Range("A28").Style = $ws.Range("A30").Style
$ws.Range("B28").Style = $ws.Range("B30").Style
$ws.Range("C28").Style = $ws.Range("C30").Style
$ws.Range("D28").Style = $ws.Range("D30").Style
$ws.Range("E28").Style = $ws.Range("E30").Style

$ws.Range("A29").Style = $ws.Range("A31").Style
$ws.Range("B29").Style = $ws.Range("B31").Style
$ws.Range("C29").Style = $ws.Range("C31").Style
$ws.Range("D29").Style = $ws.Range("D31").Style
$ws.Range("E29").Style = $ws.Range("E31").Style

# Update the sheet view
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("A28:XFD29").Select()

# Re-style header cell A1 (bold font + fill + border)
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.PatternColorIndex = -4105
$ws.Range("A1").Interior.Color = 49407
$ws.Range("A1").Interior.Pattern = -4124
